# Applies the "Add files via upload" commit: fills in the previously-blank
# Merge Sort / Heap Sort rows (and Selection Sort's time cell) across the
# three tables, adds separator rows between the tables, adds a small
# legend block below the last table, and highlights a handful of cells
# with yellow / orange fills to match the new legend colors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Cell value updates
# ---------------------------------------------------------------------

# Tabela de Tempo (rows 1-7): Selection Sort time, Merge Sort row, Heap Sort row
$ws.Range("B3").Value = "43ms"

$ws.Range("B5").Value = "755 ms"
$ws.Range("C5").Value = " 825 ms"
$ws.Range("D5").Value = "939ms"

$ws.Range("B7").Value = "461ms"
$ws.Range("C7").Value = "494ms"
$ws.Range("D7").Value = "751ms"

# Tabela de Comparacoes (rows 9-15): Selection Sort count, Merge Sort row, Heap Sort row
$ws.Range("B11").Value = 1999998

$ws.Range("B13").Value = 31017854
$ws.Range("C13").Value = 20951462
$ws.Range("D13").Value = 30345578

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1

# Tabela de Trocas (rows 17-23): Merge Sort row, Heap Sort row
$ws.Range("B21").Value = 999999
$ws.Range("C21").Value = 999999
$ws.Range("D21").Value = 999999

$ws.Range("B23").Value = 999999
$ws.Range("C23").Value = 999999
$ws.Range("D23").Value = 999999

# ---------------------------------------------------------------------
# 2) New legend block (rows 25-27)
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Legendas:"
$ws.Range("A26").Value = "Não finalizou por Tempo"
$ws.Range("A27").Value = "Código não executou"

# ---------------------------------------------------------------------
# 3) Fill colors
# ---------------------------------------------------------------------
$yellow = 10092543   # RGB(255,255,153) -> FFFF99
$orange = 10079487   # RGB(255,204,153) -> FFCC99
$dark   = 855309     # RGB(13,13,13)    -> "Black, Text 1, Lighter 5%"

# Yellow highlight on the "N/A" placeholder pairs that sit next to the
# Selection Sort / Insert Sort rows in every table.
$yellowRanges = @("C2:D2", "C4:D4", "C10:D10", "C12:D12", "C18:D18", "C20:D20")
foreach ($rng in $yellowRanges) {
    $ws.Range($rng).Interior.Color = $yellow
}

# Orange highlight on the Quick Sort row in every table.
$orangeRanges = @("B6:D6", "B14:D14", "B22:D22")
foreach ($rng in $orangeRanges) {
    $ws.Range($rng).Interior.Color = $orange
}

# Dark separator rows between tables.
$darkRanges = @("A8:D8", "A16:D16", "A24:D24")
foreach ($rng in $darkRanges) {
    $ws.Range($rng).Interior.Color = $dark
}

# Header-style formatting for the legend title (bold + centered, same
# look as the existing table headers in row 1/9/17). Alignment first,
# then bold, so this lands on the same combined style the headers use
# instead of allocating a fresh one.
$ws.Range("A25").HorizontalAlignment = -4108
$ws.Range("A25").Font.Bold = $true

# Legend swatches: "Nao finalizou por Tempo" uses the yellow highlight,
# "Codigo nao executou" uses the orange highlight. Alignment first, then
# fill, mirrors the order above so it reuses the same combined styles.
$ws.Range("A26").HorizontalAlignment = -4108
$ws.Range("A26").Interior.Color = $yellow
$ws.Range("A27").HorizontalAlignment = -4108
$ws.Range("A27").Interior.Color = $orange

# ---------------------------------------------------------------------
# 4) Selection / view state
# ---------------------------------------------------------------------
$ws.Range("B25:B26").Select()
